$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.687.37'
$ws.Range('E2').Value = '  -1.90%  '
$ws.Range('D3').Value = '2.025.82'
$ws.Range('E3').Value = '  +0.56%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '235.84'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -9.54%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.602'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.82%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '55.03'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -3.39%  '
$ws.Range('E9').Value = '  -3.58%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '58.04'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +3.46%  '
$ws.Range('E11').Value = '  -3.55%  '
$ws.Range('E12').Value = '  -0.32%  '
$ws.Range('D13').Value = '2.324.78'
$ws.Range('E13').Value = '  +0.59%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '14.20'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -1.03%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '20.25'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -6.67%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.764'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -4.32%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '5.12'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -2.13%  '
$ws.Range('D18').Value = '2.022.34'
$ws.Range('E18').Value = '  -0.33%  '
$ws.Range('D19').Value = '36.503.77'
$ws.Range('E19').Value = '  -2.22%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '67.87'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -3.27%  '
$ws.Range('D21').Value = '0.0₃0800'
$ws.Range('E21').Value = '  -4.77%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '5.38'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +5.05%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '221.44'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -5.62%  '
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.39'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.84%  '
$ws.Range('E26').Value = '  -7.80%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '163.53'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.96%  '
$ws.Range('E28').Value = '  -3.57%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.38'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +4.64%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.128'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.47%  '
$ws.Range('E31').Value = '  -3.67%  '
$ws.Range('E32').Value = '  -2.47%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.37'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -5.33%  '
$ws.Range('E34').Value = '  -5.74%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.45'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +3.09%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.26'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -5.41%  '
$ws.Range('E38').Value = '  -2.07%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.32'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.36%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.77'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +4.96%  '
$ws.Range('E41').Value = '  -2.34%  '
$ws.Range('D42').Value = '1.464.05'
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0205'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -3.49%  '
$ws.Range('B45').Value = 'TrustWalletToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '1.11'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -7.04%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '90.36'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.75%  '
$ws.Range('B47').Value = 'FTXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.08'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +34.96%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '15.37'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.69%  '
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.01'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -2.31%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.88'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -1.76%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '6.88'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -2.15%  '
